$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the single remaining data row (row 2): Gatid -> Duhat, refreshed numbers/contact
$ws.Range("A2").Value = "Duhat Elementary School"
$ws.Range("D2").Value = 123460
$ws.Range("E2").Value = "Duhat, Santa Cruz, Laguna"
$ws.Range("F2").Value = "Barangay Duhat"
$ws.Range("G2").Value = "ewan"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "09883273453"
$ws.Range("I2").Value = 1

# Remove the second data row (previously the "Oogong" entry) entirely
$ws.Rows.Item(3).Delete()

# Column widths shrink slightly to match the new (shorter) text after bestFit recalculation
$ws.Columns.Item(1).ColumnWidth = 27.333333333333332
$ws.Columns.Item(5).ColumnWidth = 29.666666666666668
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668
